$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Transistors x4 / BC337): Comments -> "Use Damo's to test"
$ws.Range("E10").Value = "Use Damo's to test"

# Row 11 (Motors x4 / H107-A03): From -> "eBay"
$ws.Range("D11").Value = "eBay"
